# This script reproduces the diff:
#  - Two new data rows are inserted into the sheet at row 512 (pushing the
#    former rows 512-533 down to become rows 514-535).
#  - The dimension grows from A1:R533 to A1:R535 automatically because the
#    used range now spans two additional rows.
#
# New row 512 (Papa / Rodeo / "1a (cosecha lavada)" ...)
# New row 513 (Papa / Rosara / "1a (cosecha)" ...)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 512; Excel shifts the
# existing rows 512:533 down to 514:535 and copies formatting from the
# row above (keeps the date style used in column D).
$ws.Range("A512:A513").EntireRow.Insert()

# ---- New row 512 ----
$ws.Cells.Item(512, 1).Value  = 5
$ws.Cells.Item(512, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(512, 3).Value  = "Maule"
$ws.Cells.Item(512, 4).Value  = 44753
$ws.Cells.Item(512, 5).Value  = 7
$ws.Cells.Item(512, 6).Value  = 100114001
$ws.Cells.Item(512, 7).Value  = "Papa"
$ws.Cells.Item(512, 8).Value  = "Rodeo"
$ws.Cells.Item(512, 9).Value  = "1a (cosecha lavada)"
$ws.Cells.Item(512, 10).Value = 1500
$ws.Cells.Item(512, 11).Value = 9000
$ws.Cells.Item(512, 12).Value = 9000
$ws.Cells.Item(512, 13).Value = 9000
$ws.Cells.Item(512, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(512, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(512, 16).Value = 360
$ws.Cells.Item(512, 17).Value = 25
$ws.Cells.Item(512, 18).Value = "Hortaliza"

# ---- New row 513 ----
$ws.Cells.Item(513, 1).Value  = 5
$ws.Cells.Item(513, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(513, 3).Value  = "Maule"
$ws.Cells.Item(513, 4).Value  = 44753
$ws.Cells.Item(513, 5).Value  = 7
$ws.Cells.Item(513, 6).Value  = 100114001
$ws.Cells.Item(513, 7).Value  = "Papa"
$ws.Cells.Item(513, 8).Value  = "Rosara"
$ws.Cells.Item(513, 9).Value  = "1a (cosecha)"
$ws.Cells.Item(513, 10).Value = 1500
$ws.Cells.Item(513, 11).Value = 6500
$ws.Cells.Item(513, 12).Value = 6500
$ws.Cells.Item(513, 13).Value = 6500
$ws.Cells.Item(513, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(513, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(513, 16).Value = 260
$ws.Cells.Item(513, 17).Value = 25
$ws.Cells.Item(513, 18).Value = "Hortaliza"

# Make sure the worksheet dimension reflects the new extent.
Write-Host "Done. UsedRange rows:" $ws.UsedRange.Rows.Count
